# "horaris funciona i afegit info extre tipus de bici"
#
# Adds a per-worker shift table (Treballador / Entrada / Sortida / Hores)
# in columns I:L of the "Napoles" and "Sants" sheets, and reshuffles the
# small stats block (Total Hores / Num treballadors / ... / Trikes) one
# column to the left (F:G -> E:G), adding a new "4 Wheeler" row and a
# "Total Hores" (shift hours) summary next to it in K:L.

$wb = $excel.ActiveWorkbook

function Update-ScheduleSheet {
    param(
        [string]$SheetName,
        $Workers,               # array of @(name, entrada, sortida, hores), starting at row 3
        [int]$StatStartRow,     # first stat row (was "Total Hores")
        [string]$TotalHoresFormula,   # e.g. "=SUM(G3:G14)"
        [int]$NumTreballadors,
        [int]$NumRutes,
        [int]$TotalPaquets,
        [int]$Trikes,
        [int]$FourWheeler,
        [string]$HoresTotalFormula    # e.g. "=SUM(L3:L15)"
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # --- 1. New header row (I2:L2), copying the existing header style ---
    $ws.Range("C2").Copy() | Out-Null
    $ws.Range("I2:L2").PasteSpecial(-4122) | Out-Null
    $ws.Range("I2").Value = "Treballador"
    $ws.Range("J2").Value = "Entrada"
    $ws.Range("K2").Value = "Sortida"
    $ws.Range("L2").Value = "Hores"

    # --- 2. Per-worker shift data (I/J/K/L), starting row 3 ---
    $r = 3
    foreach ($w in $Workers) {
        $ws.Cells.Item($r, 9).Value = $w[0]
        $ws.Cells.Item($r, 10).Value = $w[1]
        $ws.Cells.Item($r, 11).Value = $w[2]
        $ws.Cells.Item($r, 12).Value = $w[3]
        $r = $r + 1
    }

    # --- 3. Clear the old F:G stat block before rewriting it shifted to E:G ---
    $clearRange = $ws.Range($ws.Cells.Item($StatStartRow, 6), $ws.Cells.Item($StatStartRow + 5, 7))
    $clearRange.Clear() | Out-Null

    $row0 = $StatStartRow
    $ws.Range("E$row0").Value = "Total Hores"
    $ws.Range("F$row0").Formula = $TotalHoresFormula
    $ws.Range("K$row0").Value = "Total Hores"
    $ws.Range("L$row0").Formula = $HoresTotalFormula

    $row1 = $StatStartRow + 1
    $ws.Range("E$row1").Value = "Num treballadors"
    $ws.Range("F$row1").Value = $NumTreballadors

    $row2 = $StatStartRow + 2
    $ws.Range("E$row2").Value = "Num Rutes"
    $ws.Range("F$row2").Value = $NumRutes

    $row3 = $StatStartRow + 3
    $ws.Range("E$row3").Value = "Total Paquets"
    $ws.Range("F$row3").Value = $TotalPaquets

    $row4 = $StatStartRow + 4
    $ws.Range("E$row4").Value = "Trikes"
    $ws.Range("F$row4").Value = $Trikes
    $ws.Range("G$row4").Formula = "=G$row4/G$row2"

    $row5 = $StatStartRow + 5
    $ws.Range("E$row5").Value = "4 Wheeler"
    $ws.Range("F$row5").Value = $FourWheeler
    $ws.Range("G$row5").Formula = "=G$row5/G$row2"
}

# ----- Napoles -----
$napolesWorkers = @(
    ,@("Zoe", "7:45", "13:45", 6)
    ,@("Fede Goss", "8:15", "12:15", 4)
    ,@("Gianluca", "8:30", "12:30", 4)
    ,@("Laila", "9:00", "11:30", 2.5)
    ,@("Vladi", "9:15", "16:15", 7)
    ,@("Erick", "9:45", "12:15", 2.5)
    ,@("Sebastián", "10:00", "13:00", 3)
)
Update-ScheduleSheet "Napoles" $napolesWorkers 16 "=SUM(G3:G14)" 12 19 77 8 11 "=SUM(L3:L15)"

# ----- Sants -----
$santsWorkers = @(
    ,@("Jordi", "7:30", "12:30", 5)
    ,@("Rocco", "7:45", "12:00", 4.3)
    ,@("Alejandro", "8:00", "15:00", 7)
    ,@("Hans", "9:30", "12:00", 2.5)
    ,@("Jaime", "10:00", "14:00", 4)
    ,@("Diego", "16:00", "20:00", 4)
)
Update-ScheduleSheet "Sants" $santsWorkers 11 "=SUM(G3:G9)" 7 19 79 10 9 "=SUM(L3:L10)"

Write-Host "Done applying schedule updates"
